$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new rows 4 and 6 (same pattern as row 2) and row 5 (same pattern as row 3) ---
# Copying the whole row preserves per-cell styles (date format on H, hyperlink style on J)
# without introducing redundant style entries.
$ws.Range("A2:Z2").Copy($ws.Range("A4:Z4"))
$ws.Range("A3:Z3").Copy($ws.Range("A5:Z5"))
$ws.Range("A2:Z2").Copy($ws.Range("A6:Z6"))

# --- ids / dni for the new rows ---
$ws.Range("A4").Value = 16200277
$ws.Range("B4").Value = 71395616
$ws.Range("A5").Value = 16200278
$ws.Range("B5").Value = 71395617
$ws.Range("A6").Value = 16200279
$ws.Range("B6").Value = 71395618

# --- alumnoPlan (K) for the new rows ---
$ws.Range("K4").Value = 17
$ws.Range("K5").Value = 18
$ws.Range("K6").Value = 19

# --- ocCorreoPersonal (Z) for the new rows ---
$ws.Range("Z4").Value = "aaaaa"
$ws.Range("Z5").Value = "bbbb"
$ws.Range("Z6").Value = "jsjsjs"

# --- alumnoCorreoIns (J) text for every row with a hyperlink ---
$ws.Range("J2").Value = "a@unmsm.edu.pe"
$ws.Range("J3").Value = "b@unmsm.edu.pe"
$ws.Range("J4").Value = "c@unmsm.edu.pe"
$ws.Range("J5").Value = "d@unmsm.edu.pe"
$ws.Range("J6").Value = "e@unmsm.edu.pe"

# --- Rebuild the mailto hyperlinks for column J (order matches target relationship ids) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:a@unmsm.edu.pe") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J3"), "mailto:b@unmsm.edu.pe") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J4"), "mailto:c@unmsm.edu.pe") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J6"), "mailto:e@unmsm.edu.pe") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J5"), "mailto:d@unmsm.edu.pe") | Out-Null

# Hyperlinks.Add() re-stamps the affected cells with a fresh style; restore the
# workbook's existing "Hipervínculo" cell style so no cell style drifts.
$ws.Range("J2:J6").Style = "Hipervínculo"

# --- View state: scroll so column L is the leftmost visible column, then land on Z6 ---
$ws.Range("L1").Select()
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("Z6").Select()
